# Bug fixes in fixtures update script:
# Populate the village / mediator / group fixture sheets with data rows,
# and fix the mediator/group header (C1 should read "field: village_id").
#
# Numeric-looking ids ("17583", "6003", ...) must be written as TEXT
# (not numbers) to match the lookup-fixture convention used elsewhere in
# this workbook. We force text entry with a leading apostrophe and then
# reset the cell style back to Normal so no stray number-format sticks to
# the cell.

function Set-TextValue {
    param($Cell, [string]$Text)

    if ($Text -match '^-?[0-9]+(\.[0-9]+)?$') {
        # Looks like a number to Excel's auto-detection - force text entry.
        $Cell.Value = "'" + $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# village sheet
# ---------------------------------------------------------------------
$village = $wb.Worksheets.Item("village")

$villageRows = @(
    @("17583", "Patimatla", "abezu"),
    @("18932", "Sangam (Simra)", "abezu")
)

$r = 2
foreach ($row in $villageRows) {
    Set-TextValue $village.Cells.Item($r, 1) $row[0]
    Set-TextValue $village.Cells.Item($r, 2) $row[1]
    Set-TextValue $village.Cells.Item($r, 3) $row[2]
    $r++
}

# ---------------------------------------------------------------------
# mediator sheet
# ---------------------------------------------------------------------
$mediator = $wb.Worksheets.Item("mediator")

# header fix: C1 should read "field: village_id"
Set-TextValue $mediator.Cells.Item(1, 3) "field: village_id"

$mediatorRows = @(
    @("6003", "M.Narsaiah (CA)", "17583", "abezu"),
    @("7906", "Sohan Sah", "18932", "abezu")
)

$r = 2
foreach ($row in $mediatorRows) {
    Set-TextValue $mediator.Cells.Item($r, 1) $row[0]
    Set-TextValue $mediator.Cells.Item($r, 2) $row[1]
    Set-TextValue $mediator.Cells.Item($r, 3) $row[2]
    Set-TextValue $mediator.Cells.Item($r, 4) $row[3]
    $r++
}

# ---------------------------------------------------------------------
# group sheet
# ---------------------------------------------------------------------
$group = $wb.Worksheets.Item("group")

# header fix: C1 should read "field: village_id"
Set-TextValue $group.Cells.Item(1, 3) "field: village_id"

$groupRows = @(
    @("70844", "Adarsha Mahila", "17583", "abezu"),
    @("70845", "Bharathmatha", "17583", "abezu"),
    @("70846", "Deevena", "17583", "abezu"),
    @("70847", "Jhansi", "17583", "abezu"),
    @("70848", "Kanakadurga", "17583", "abezu"),
    @("70849", "Mahalaxmi", "17583", "abezu"),
    @("70850", "Mallikarjuna", "17583", "abezu"),
    @("70851", "Thelugu Mahila", "17583", "abezu"),
    @("70852", "Sri Laxmi Durga", "17583", "abezu"),
    @("70853", "Sri Rama", "17583", "abezu"),
    @("78952", "Sharda SHG", "18932", "abezu"),
    @("78953", "Anju", "18932", "abezu"),
    @("78954", "Eklakhiya", "18932", "abezu"),
    @("78955", "Gausiya", "18932", "abezu"),
    @("78956", "Sita", "18932", "abezu"),
    @("78957", "Madarsa", "18932", "abezu"),
    @("78958", "Bhawani", "18932", "abezu"),
    @("78959", "Jyoti", "18932", "abezu"),
    @("78960", "Chameli", "18932", "abezu"),
    @("78961", "Tegiya", "18932", "abezu"),
    @("78962", "Radha", "18932", "abezu")
)

$r = 2
foreach ($row in $groupRows) {
    Set-TextValue $group.Cells.Item($r, 1) $row[0]
    Set-TextValue $group.Cells.Item($r, 2) $row[1]
    Set-TextValue $group.Cells.Item($r, 3) $row[2]
    Set-TextValue $group.Cells.Item($r, 4) $row[3]
    $r++
}
